$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.193.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -5.60%  '
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.205.94'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -8.58%  '
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.15'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.39%  '
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.52'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -13.08%  '
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.206.44'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -8.53%  '
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -13.10%  '
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -12.52%  '
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -18.98%  '
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -19.65%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.25'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -19.34%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000231'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -15.79%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.714.10'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -8.78%  '
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.222.59'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -8.01%  '
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.994.69'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.98%  '
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.115'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.53%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '513.08'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -15.88%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.82'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -17.68%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.43'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -16.73%  '
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.742'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -15.28%  '
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.42'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -17.45%  '
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.77'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -14.96%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.82'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -17.44%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.984'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.22'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -13.10%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '28.19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -16.36%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.04'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -20.09%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.21'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -19.50%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Mantle'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -13.16%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.43'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -17.72%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '516.43'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -19.48%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -19.84%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.18'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -22.90%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.80'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.06%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0412'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -12.06%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.99'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -16.02%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0821'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -17.05%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.110'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -22.20%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.758.09'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -17.90%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.55'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -28.39%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -16.88%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.246'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -19.78%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.68'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.27%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₃0530'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -28.27%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.70'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -23.15%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.110'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -14.77%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Fetch.AI'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -21.22%  '
